$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge "palette" + "Pico" runs into a single "palettePico" run.
#    (Replacing the text with itself forces Word to collapse the matched
#    span into one run.)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("palettePico", $true, $false, $false, $false, $false, `
    $true, 1, $false, "palettePico", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Merge "try " + "here " + "> " runs into a single run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("try here > ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "try here > ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Merge the "graphics automatically ... into your game." runs into one.
#    The replace also swallows the run boundary that precedes "graphics" (it
#    shares the same formatting), so re-split the run there afterwards to
#    keep the preceding "...compiles " run intact, matching the source.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute( `
    "graphics automatically into data files that can then be directly incorporated into your game.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "graphics automatically into data files that can then be directly incorporated into your game.", `
    2) | Out-Null

$anchor = $d.Content
$anchor.Find.Execute("compiles graphics automatically into data files that can then be directly incorporated into your game.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$sub = $d.Range($anchor.Start, $anchor.End)
$sub.Find.Execute("graphics automatically into data files that can then be directly incorporated into your game.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("zzTempSplit0", $sub)
$d.Bookmarks("zzTempSplit0").Delete()

# ---------------------------------------------------------------------------
# 4) Remove the old _GoBack bookmark (next to "Program").
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------------
# 5) "byte" -> "uint8_t" in the "const byte puzzles_5x5_count = 2;" sample,
#    keeping it split across three runs like Word would after an in-place
#    edit: " " / "uint8_t" / " puzzles_5x5_count = 2;". The leading space is
#    part of the match so we don't touch the (proofErr-separated) "const"
#    run that precedes it.
# ---------------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute(" byte puzzles_5x5_count = 2;", $true, $false, $false, $false, `
    $false, $true, 1, $false, " uint8_t puzzles_5x5_count = 2;", 2) | Out-Null

$anchor = $d.Content
$anchor.Find.Execute(" uint8_t puzzles_5x5_count = 2;", $true, $false, $false, $false, `
    $false, $true, 1, $false, "", 0) | Out-Null
$sub = $d.Range($anchor.Start, $anchor.End)
$sub.Find.Execute("uint8_t", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("zzTempSplit1", $sub)
$d.Bookmarks("zzTempSplit1").Delete()

# ---------------------------------------------------------------------------
# 6) Expand "(bytes, integers or even other objects)" into
#    "- bytes (known as uint8_t in C++), integers or even other objects"
#    split across four runs as in the diff.
# ---------------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute("(bytes, integers or even other objects)", $true, $false, $false, `
    $false, $false, $true, 1, $false, `
    "- bytes (known as uint8_t in C++), integers or even other objects", 2) | Out-Null

$anchor = $d.Content
$anchor.Find.Execute("- bytes (known as uint8_t in C++), integers or even other objects", $true, `
    $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$scope = $d.Range($anchor.Start, $anchor.End)

$sub = $d.Range($scope.Start, $scope.End)
$sub.Find.Execute("- ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("zzTempSplit2", $sub)
$d.Bookmarks("zzTempSplit2").Delete()

$sub = $d.Range($scope.Start, $scope.End)
$sub.Find.Execute("bytes", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("zzTempSplit3", $sub)
$d.Bookmarks("zzTempSplit3").Delete()

$sub = $d.Range($scope.Start, $scope.End)
$sub.Find.Execute(" (known as uint8_t in C++)", $true, $false, $false, $false, $false, $true, 1, `
    $false, "", 0) | Out-Null
$d.Bookmarks.Add("zzTempSplit4", $sub)
$d.Bookmarks("zzTempSplit4").Delete()

# ---------------------------------------------------------------------------
# 7) Move the _GoBack bookmark onto the "byte" in the hexadecimal-notation
#    paragraph, splitting that run into three pieces around it.
# ---------------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute("combining two values into one byte.  When we read", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0) | Out-Null

$sub = $d.Range($anchor.Start, $anchor.End)
$sub.Find.Execute("byte", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("_GoBack", $sub)

Write-Output "done"
